# Updates a set of numeric values inside the single results table of the
# document. Cells are addressed by (row, column) within the table (the
# header row is row 1), and each cell's current text is verified against
# the expected "old" value before being overwritten, to avoid silently
# touching the wrong cell.
#
# NOTE: a `TableCell.Range` object, once a *different* part of the
# document has already been edited, can no longer be trusted as the
# scope for a subsequent `.Find.Execute(...)` call on this runtime - the
# search ends up running against stale offsets and can land on (and
# corrupt) unrelated earlier text. Re-deriving a brand-new `$d.Range(start,
# end)` from the cell's current Start/End right before each Find call
# keeps every substitution confined to the intended cell.

$d = $word.ActiveDocument

function Set-CellValue($Row, $Col, $OldVal, $NewVal) {
    $table = $d.Tables.Item(1)
    $cell = $table.Cell($Row, $Col)
    $cellRange = $cell.Range
    $range = $d.Range($cellRange.Start, $cellRange.End)

    # Cell ranges include the trailing cell-mark / paragraph-mark
    # characters, so trim those off before comparing / writing text.
    $text = $range.Text
    $text = $text.TrimEnd([char]7, [char]13, [char]10)

    if ($text -ne $OldVal) {
        throw "Cell ($Row,$Col) expected '$OldVal' but found '$text'"
    }

    $find = $range.Find
    $find.ClearFormatting()
    # MatchCase:=true, Wrap:=wdFindStop(0), Replace:=wdReplaceOne(1) so the
    # substitution stays confined to this single table cell instead of
    # leaking into sibling cells that happen to share the same old value.
    $find.Execute($OldVal, $true, $false, $false, $false, $false, $true, 0, $false, $NewVal, 1)
}

Set-CellValue 2 8 "-3.197" "-3.196"
Set-CellValue 2 11 "6843.586" "2767.478"
Set-CellValue 2 13 "100" "98"
Set-CellValue 3 8 "-2.719" "-2.718"
Set-CellValue 3 11 "6811.205" "2735.791"
Set-CellValue 3 13 "100" "98"
Set-CellValue 4 8 "-2.840" "-2.839"
Set-CellValue 4 11 "6453.431" "2391.005"
Set-CellValue 4 13 "100" "99"
Set-CellValue 5 8 "2.400" "2.401"
Set-CellValue 5 11 "1636.931" "1636.793"
Set-CellValue 6 8 "2.372" "2.373"
Set-CellValue 6 11 "1605.365" "1605.231"
Set-CellValue 7 8 "1.648" "1.649"
Set-CellValue 7 11 "1502.464" "1502.321"
Set-CellValue 8 8 "-0.945" "-0.944"
Set-CellValue 8 11 "96.206" "96.205"
Set-CellValue 9 8 "-0.945" "-0.944"
Set-CellValue 9 11 "96.206" "96.205"
Set-CellValue 10 8 "-5.626" "-5.625"
Set-CellValue 10 11 "59.137" "59.146"
Set-CellValue 11 11 "16.566" "16.585"
Set-CellValue 11 12 "0.280" "0.279"
Set-CellValue 12 11 "16.566" "16.585"
Set-CellValue 12 12 "0.280" "0.279"
Set-CellValue 13 11 "8.438" "8.458"
Set-CellValue 13 12 "0.586" "0.584"
Set-CellValue 14 8 "-2.097" "-2.092"
Set-CellValue 14 11 "307.160" "200.514"
Set-CellValue 14 13 "76" "68"
Set-CellValue 15 8 "-2.097" "-2.092"
Set-CellValue 15 11 "307.160" "200.514"
Set-CellValue 15 13 "76" "68"
Set-CellValue 16 6 "-0.138" "-0.137"
Set-CellValue 16 7 "0.032" "0.031"
Set-CellValue 16 8 "-4.348" "-4.358"
Set-CellValue 16 11 "75.098" "72.330"
Set-CellValue 16 13 "30" "26"
Set-CellValue 17 8 "-2.790" "-2.789"
Set-CellValue 17 11 "3.763" "3.782"
Set-CellValue 18 8 "-2.790" "-2.789"
Set-CellValue 18 11 "3.763" "3.782"
Set-CellValue 19 11 "3.270" "3.289"
Set-CellValue 19 12 "0.987" "0.986"

Write-Output "Applied 44 cell updates successfully."
